$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155, shifting the existing rows 155-177 down
# to 156-178 (this is the "weekly" data point being added, per the commit
# message "Fruta / hortaliza, semanal").
$ws.Rows.Item(155).Insert()

$ws.Cells.Item(155, 1).Value = 7
$ws.Cells.Item(155, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(155, 3).Value = "Ñuble"
$ws.Cells.Item(155, 4).Value = 44491
$ws.Cells.Item(155, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(155, 5).Value = 16
$ws.Cells.Item(155, 6).Value = 100112023
$ws.Cells.Item(155, 7).Value = "Brócoli"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 240
$ws.Cells.Item(155, 11).Value = 750
$ws.Cells.Item(155, 12).Value = 800
$ws.Cells.Item(155, 13).Value = 775
$ws.Cells.Item(155, 14).Value = "`$/unidad"
$ws.Cells.Item(155, 15).Value = "Región Metropolitana"
$ws.Cells.Item(155, 16).Value = 775
$ws.Cells.Item(155, 17).Value = 1
$ws.Cells.Item(155, 18).Value = "Hortaliza"
